$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 192, shifting existing rows 192:256 down to 193:257
$ws.Rows("192:192").Insert()

# Populate the newly inserted row 192 with the new weekly record
$ws.Range("A192").Value2 = 11
$ws.Range("B192").Value2 = "Vega Monumental Concepción"
$ws.Range("C192").Value2 = "Bíobío"
$ws.Range("D192").Value2 = 44588
$ws.Range("E192").Value2 = 8
$ws.Range("F192").Value2 = 100114014
$ws.Range("G192").Value2 = "Betarraga"
$ws.Range("H192").Value2 = "Sin especificar"
$ws.Range("I192").Value2 = "Primera"
$ws.Range("J192").Value2 = 500
$ws.Range("K192").Value2 = 600
$ws.Range("L192").Value2 = 650
$ws.Range("M192").Value2 = 620
$ws.Range("N192").Value2 = "`$/paquete 5 unidades"
$ws.Range("O192").Value2 = "Región Metropolitana"
$ws.Range("P192").Value2 = 124
$ws.Range("Q192").Value2 = 5
$ws.Range("R192").Value2 = "Hortaliza"
